# Depression and Poverty.pptx - week 4 refresh
#
# Updates the title-slide subtitle box ("Subtitle 2" / shape id 3 on slide 1):
#   - the date line "October 16, 2022"            -> "October 23, 2022"
#   - the repo-link text "...ANA500-Week-3"        -> "...ANA500-Week-4"
#   - removes the now-unused trailing blank paragraph under the link

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the subtitle placeholder shape by name rather than a hard-coded index.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Subtitle 2") {
        $shp = $candidate
        break
    }
}

if ($shp -ne $null) {
    $tr = $shp.TextFrame.TextRange

    # Update the date run in place (keeps its existing run formatting).
    $foundDate = $tr.Find("October 16, 2022", 0)
    if ($foundDate -ne $null) {
        $foundDate.Text = "October 23, 2022"
    }

    # Update the hyperlink run's visible text (exact, whole-run match keeps
    # the single run / hlinkClick intact rather than splitting it).
    $foundUrl = $tr.Find("https://github.com/danlagos/ANA500-Week-3", 0)
    if ($foundUrl -ne $null) {
        $foundUrl.Text = "https://github.com/danlagos/ANA500-Week-4"
    }

    # Drop the trailing empty paragraph that used to sit below the link.
    $paraCount = $tr.Paragraphs(0, -1).Count
    $lastPara = $tr.Paragraphs($paraCount + 1, 1)
    if ($lastPara -ne $null -and $lastPara.Length -eq 0) {
        $lastPara.Delete()
    }
}
